# Fix typo "Transaksi Berahasil Dibuat" -> "Transaksi Berhasil Dibuat"
# on the "Order Tracking" sheet, then leave that sheet active with the
# selection on J10 (mirroring the author's final UI state).

$wb = $excel.ActiveWorkbook

$tracking = $wb.Worksheets.Item("Order Tracking")

$fixedRows = @(2, 3, 4, 5, 7, 8, 9, 11, 13, 14)
foreach ($r in $fixedRows) {
    $cell = $tracking.Cells.Item($r, 4)
    if ($cell.Value() -eq "Transaksi Berahasil Dibuat") {
        $cell.Value = "Transaksi Berhasil Dibuat"
    }
}

$tracking.Activate()
$tracking.Range("J10").Select()
